$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "model_8_2_0"
$ws.Cells.Item(2, 2).Value = 0.6501110450060772
$ws.Cells.Item(2, 3).Value = -0.2284759110106034
$ws.Cells.Item(2, 4).Value = -0.2869403027450084
$ws.Cells.Item(2, 5).Value = -0.05205087385175955
$ws.Cells.Item(2, 6).Value = 0.3872239887714386
$ws.Cells.Item(2, 7).Value = 1.441739320755005
$ws.Cells.Item(2, 8).Value = 0.6524465680122375
$ws.Cells.Item(2, 9).Value = 1.070307731628418

$ws.Cells.Item(3, 1).Value = "model_8_2_1"
$ws.Cells.Item(3, 2).Value = 0.6684568114876195
$ws.Cells.Item(3, 3).Value = -0.1978381252582904
$ws.Cells.Item(3, 4).Value = -0.4178504393717983
$ws.Cells.Item(3, 5).Value = -0.06403831400476112
$ws.Cells.Item(3, 6).Value = 0.366920679807663
$ws.Cells.Item(3, 7).Value = 1.40578281879425
$ws.Cells.Item(3, 8).Value = 0.7188147306442261
$ws.Cells.Item(3, 9).Value = 1.082503199577332

$ws.Cells.Item(4, 1).Value = "model_8_2_3"
$ws.Cells.Item(4, 2).Value = 0.691642624404355
$ws.Cells.Item(4, 3).Value = 0.01192194045987316
$ws.Cells.Item(4, 4).Value = -1.034139813351757
$ws.Cells.Item(4, 5).Value = -0.08045851275184779
$ws.Cells.Item(4, 6).Value = 0.3412607312202454
$ws.Cells.Item(4, 7).Value = 1.159608364105225
$ws.Cells.Item(4, 8).Value = 1.0312579870224
$ws.Cells.Item(4, 9).Value = 1.099208235740662

$ws.Cells.Item(5, 1).Value = "model_8_2_4"
$ws.Cells.Item(5, 2).Value = 0.6937162513163642
$ws.Cells.Item(5, 3).Value = 0.020136194692009
$ws.Cells.Item(5, 4).Value = -1.040765842987379
$ws.Cells.Item(5, 5).Value = -0.07699364442631551
$ws.Cells.Item(5, 6).Value = 0.3389658629894257
$ws.Cells.Item(5, 7).Value = 1.149968147277832
$ws.Cells.Item(5, 8).Value = 1.034617185592651
$ws.Cells.Item(5, 9).Value = 1.095683336257935

$ws.Cells.Item(6, 1).Value = "model_8_2_6"
$ws.Cells.Item(6, 2).Value = 0.6966624769836658
$ws.Cells.Item(6, 3).Value = 0.03124661477850155
$ws.Cells.Item(6, 4).Value = -1.04546742598412
$ws.Cells.Item(6, 5).Value = -0.07131235925591106
$ws.Cells.Item(6, 6).Value = 0.3357052505016327
$ws.Cells.Item(6, 7).Value = 1.136929035186768
$ws.Cells.Item(6, 8).Value = 1.037000894546509
$ws.Cells.Item(6, 9).Value = 1.089903473854065

$ws.Cells.Item(7, 1).Value = "model_8_2_5"
$ws.Cells.Item(7, 2).Value = 0.6967424061097212
$ws.Cells.Item(7, 3).Value = 0.0325015632072132
$ws.Cells.Item(7, 4).Value = -1.039663732075012
$ws.Cells.Item(7, 5).Value = -0.06918369787986389
$ws.Cells.Item(7, 6).Value = 0.335616797208786
$ws.Cells.Item(7, 7).Value = 1.135456204414368
$ws.Cells.Item(7, 8).Value = 1.034058570861816
$ws.Cells.Item(7, 9).Value = 1.087737798690796

$ws.Cells.Item(8, 1).Value = "model_8_2_7"
$ws.Cells.Item(8, 2).Value = 0.6968616334442568
$ws.Cells.Item(8, 3).Value = 0.03038969304723527
$ws.Cells.Item(8, 4).Value = -1.044859842917171
$ws.Cells.Item(8, 5).Value = -0.07169477239891564
$ws.Cells.Item(8, 6).Value = 0.3354848325252533
$ws.Cells.Item(8, 7).Value = 1.137934684753418
$ws.Cells.Item(8, 8).Value = 1.03669273853302
$ws.Cells.Item(8, 9).Value = 1.090292572975159

$ws.Cells.Item(9, 1).Value = "model_8_2_8"
$ws.Cells.Item(9, 2).Value = 0.6975552720109063
$ws.Cells.Item(9, 3).Value = 0.0323578303608687
$ws.Cells.Item(9, 4).Value = -1.050848659991428
$ws.Cells.Item(9, 5).Value = -0.07189680784663599
$ws.Cells.Item(9, 6).Value = 0.3347172141075134
$ws.Cells.Item(9, 7).Value = 1.135624766349792
$ws.Cells.Item(9, 8).Value = 1.039728879928589
$ws.Cells.Item(9, 9).Value = 1.090498089790344

$ws.Cells.Item(10, 1).Value = "model_8_2_10"
$ws.Cells.Item(10, 2).Value = 0.6996089075163761
$ws.Cells.Item(10, 3).Value = 0.03878159183018104
$ws.Cells.Item(10, 4).Value = -1.043675682176504
$ws.Cells.Item(10, 5).Value = -0.06629152750899192
$ws.Cells.Item(10, 6).Value = 0.332444429397583
$ws.Cells.Item(10, 7).Value = 1.128085851669312
$ws.Cells.Item(10, 8).Value = 1.036092519760132
$ws.Cells.Item(10, 9).Value = 1.084795594215393

$ws.Cells.Item(11, 1).Value = "model_8_2_9"
$ws.Cells.Item(11, 2).Value = 0.699814601877117
$ws.Cells.Item(11, 3).Value = 0.03918054699337492
$ws.Cells.Item(11, 4).Value = -1.039085181707728
$ws.Cells.Item(11, 5).Value = -0.06497032077728071
$ws.Cells.Item(11, 6).Value = 0.3322167992591858
$ws.Cells.Item(11, 7).Value = 1.127617716789246
$ws.Cells.Item(11, 8).Value = 1.033765077590942
$ws.Cells.Item(11, 9).Value = 1.083451271057129

$ws.Cells.Item(12, 1).Value = "model_8_2_11"
$ws.Cells.Item(12, 2).Value = 0.7031287059214336
$ws.Cells.Item(12, 3).Value = 0.04205024970297155
$ws.Cells.Item(12, 4).Value = -1.011910822952106
$ws.Cells.Item(12, 5).Value = -0.05684837325465453
$ws.Cells.Item(12, 6).Value = 0.3285490572452545
$ws.Cells.Item(12, 7).Value = 1.124249935150146
$ws.Cells.Item(12, 8).Value = 1.019988536834717
$ws.Cells.Item(12, 9).Value = 1.075188517570496

$ws.Cells.Item(13, 1).Value = "model_8_2_12"
$ws.Cells.Item(13, 2).Value = 0.7035206987778346
$ws.Cells.Item(13, 3).Value = 0.04293746120300779
$ws.Cells.Item(13, 4).Value = -1.009252619927862
$ws.Cells.Item(13, 5).Value = -0.05568230346306713
$ws.Cells.Item(13, 6).Value = 0.3281152248382568
$ws.Cells.Item(13, 7).Value = 1.123208522796631
$ws.Cells.Item(13, 8).Value = 1.018640756607056
$ws.Cells.Item(13, 9).Value = 1.074002146720886

$ws.Cells.Item(14, 1).Value = "model_8_2_13"
$ws.Cells.Item(14, 2).Value = 0.7037712431225616
$ws.Cells.Item(14, 3).Value = 0.04390346390768318
$ws.Cells.Item(14, 4).Value = -1.008551266305151
$ws.Cells.Item(14, 5).Value = -0.05492494172548756
$ws.Cells.Item(14, 6).Value = 0.3278379440307617
$ws.Cells.Item(14, 7).Value = 1.122074842453003
$ws.Cells.Item(14, 8).Value = 1.018285274505615
$ws.Cells.Item(14, 9).Value = 1.07323169708252

$ws.Cells.Item(15, 1).Value = "model_8_2_14"
$ws.Cells.Item(15, 2).Value = 0.7039983757709021
$ws.Cells.Item(15, 3).Value = 0.04500844951477445
$ws.Cells.Item(15, 4).Value = -1.008391969522505
$ws.Cells.Item(15, 5).Value = -0.05421013865179836
$ws.Cells.Item(15, 6).Value = 0.3275865912437439
$ws.Cells.Item(15, 7).Value = 1.12077808380127
$ws.Cells.Item(15, 8).Value = 1.0182044506073
$ws.Cells.Item(15, 9).Value = 1.07250452041626

$ws.Cells.Item(16, 1).Value = "model_8_2_15"
$ws.Cells.Item(16, 2).Value = 0.7042327654239936
$ws.Cells.Item(16, 3).Value = 0.04621969836469719
$ws.Cells.Item(16, 4).Value = -1.00828130329143
$ws.Cells.Item(16, 5).Value = -0.05344758010656414
$ws.Cells.Item(16, 6).Value = 0.3273271918296814
$ws.Cells.Item(16, 7).Value = 1.119356632232666
$ws.Cells.Item(16, 8).Value = 1.018148422241211
$ws.Cells.Item(16, 9).Value = 1.071728706359863

$ws.Cells.Item(17, 1).Value = "model_8_2_16"
$ws.Cells.Item(17, 2).Value = 0.704478071791208
$ws.Cells.Item(17, 3).Value = 0.04751543107822664
$ws.Cells.Item(17, 4).Value = -1.008179240135069
$ws.Cells.Item(17, 5).Value = -0.05263378690501441
$ws.Cells.Item(17, 6).Value = 0.3270556926727295
$ws.Cells.Item(17, 7).Value = 1.117835879325867
$ws.Cells.Item(17, 8).Value = 1.018096685409546
$ws.Cells.Item(17, 9).Value = 1.070900797843933

$ws.Cells.Item(18, 1).Value = "model_8_2_17"
$ws.Cells.Item(18, 2).Value = 0.7047355441063143
$ws.Cells.Item(18, 3).Value = 0.04887759477256337
$ws.Cells.Item(18, 4).Value = -1.008051141808435
$ws.Cells.Item(18, 5).Value = -0.05176879775691412
$ws.Cells.Item(18, 6).Value = 0.3267707526683807
$ws.Cells.Item(18, 7).Value = 1.116237163543701
$ws.Cells.Item(18, 8).Value = 1.018031716346741
$ws.Cells.Item(18, 9).Value = 1.070020794868469

$ws.Cells.Item(19, 1).Value = "model_8_2_18"
$ws.Cells.Item(19, 2).Value = 0.7050008187282305
$ws.Cells.Item(19, 3).Value = 0.05029031934991191
$ws.Cells.Item(19, 4).Value = -1.007894644046823
$ws.Cells.Item(19, 5).Value = -0.05087117656134588
$ws.Cells.Item(19, 6).Value = 0.3264771401882172
$ws.Cells.Item(19, 7).Value = 1.114579319953918
$ws.Cells.Item(19, 8).Value = 1.0179523229599
$ws.Cells.Item(19, 9).Value = 1.06910765171051

$ws.Cells.Item(20, 1).Value = "model_8_2_19"
$ws.Cells.Item(20, 2).Value = 0.7050678249318738
$ws.Cells.Item(20, 3).Value = 0.05064912868551275
$ws.Cells.Item(20, 4).Value = -1.007856261778156
$ws.Cells.Item(20, 5).Value = -0.05064404158262836
$ws.Cells.Item(20, 6).Value = 0.3264029920101166
$ws.Cells.Item(20, 7).Value = 1.114158153533936
$ws.Cells.Item(20, 8).Value = 1.017932891845703
$ws.Cells.Item(20, 9).Value = 1.068876624107361

$ws.Cells.Item(21, 1).Value = "model_8_2_20"
$ws.Cells.Item(21, 2).Value = 0.705974824941326
$ws.Cells.Item(21, 3).Value = 0.05377372365305899
$ws.Cells.Item(21, 4).Value = -1.00278834075388
$ws.Cells.Item(21, 5).Value = -0.04754644329157265
$ws.Cells.Item(21, 6).Value = 0.3253992199897766
$ws.Cells.Item(21, 7).Value = 1.110491275787354
$ws.Cells.Item(21, 8).Value = 1.015363574028015
$ws.Cells.Item(21, 9).Value = 1.065725088119507

$ws.Cells.Item(22, 1).Value = "model_8_2_21"
$ws.Cells.Item(22, 2).Value = 0.7061026026742423
$ws.Cells.Item(22, 3).Value = 0.05427571074660342
$ws.Cells.Item(22, 4).Value = -1.002174491366491
$ws.Cells.Item(22, 5).Value = -0.0470965421061198
$ws.Cells.Item(22, 6).Value = 0.325257807970047
$ws.Cells.Item(22, 7).Value = 1.109902024269104
$ws.Cells.Item(22, 8).Value = 1.015052437782288
$ws.Cells.Item(22, 9).Value = 1.065267443656921

$ws.Cells.Item(23, 1).Value = "model_8_2_22"
$ws.Cells.Item(23, 2).Value = 0.7062343741146515
$ws.Cells.Item(23, 3).Value = 0.05479821376583538
$ws.Cells.Item(23, 4).Value = -1.001576144690069
$ws.Cells.Item(23, 5).Value = -0.04663972600690114
$ws.Cells.Item(23, 6).Value = 0.3251119554042816
$ws.Cells.Item(23, 7).Value = 1.109288811683655
$ws.Cells.Item(23, 8).Value = 1.014749050140381
$ws.Cells.Item(23, 9).Value = 1.064802646636963

$ws.Cells.Item(24, 1).Value = "model_8_2_23"
$ws.Cells.Item(24, 2).Value = 0.7063650299443551
$ws.Cells.Item(24, 3).Value = 0.05537683076841449
$ws.Cells.Item(24, 4).Value = -1.001238488438241
$ws.Cells.Item(24, 5).Value = -0.04620082644380696
$ws.Cells.Item(24, 6).Value = 0.3249673843383789
$ws.Cells.Item(24, 7).Value = 1.108609676361084
$ws.Cells.Item(24, 8).Value = 1.014577865600586
$ws.Cells.Item(24, 9).Value = 1.064356207847595

$ws.Cells.Item(25, 1).Value = "model_8_2_24"
$ws.Cells.Item(25, 2).Value = 0.7064618117317918
$ws.Cells.Item(25, 3).Value = 0.0558084677248083
$ws.Cells.Item(25, 4).Value = -1.000898655371177
$ws.Cells.Item(25, 5).Value = -0.04585902341623083
$ws.Cells.Item(25, 6).Value = 0.3248602747917175
$ws.Cells.Item(25, 7).Value = 1.108103275299072
$ws.Cells.Item(25, 8).Value = 1.014405608177185
$ws.Cells.Item(25, 9).Value = 1.064008474349976

$ws.Cells.Item(26, 1).Value = "model_8_2_2"
$ws.Cells.Item(26, 2).Value = 0.7086049953795519
$ws.Cells.Item(26, 3).Value = 0.01761939923708389
$ws.Cells.Item(26, 4).Value = -0.6739086711214592
$ws.Cells.Item(26, 5).Value = 0.007500468787655534
$ws.Cells.Item(26, 6).Value = 0.322488397359848
$ws.Cells.Item(26, 7).Value = 1.152921915054321
$ws.Cells.Item(26, 8).Value = 0.8486298322677612
$ws.Cells.Item(26, 9).Value = 1.009722948074341

